$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1, index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 3146
$ws1.Range("F5").Value = 2160
$ws1.Range("F8").Value = 926
$ws1.Range("F9").Value = 998
$ws1.Range("F10").Value = 238
$ws1.Range("F11").Value = 458
$ws1.Range("F16").Value = 7721
$ws1.Range("F17").Value = 337
$ws1.Range("F18").Value = 2454
$ws1.Range("F20").Value = 225
$ws1.Range("F23").Value = 530
$ws1.Range("F26").Value = 977
$ws1.Range("F28").Value = 1649
$ws1.Range("F30").Value = 1165
$ws1.Range("F34").Value = 164
$ws1.Range("F37").Value = 162
$ws1.Range("F38").Value = 333
$ws1.Range("F40").Value = 213

# Sheet "全部类型" (sheet4, index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 3146
$ws4.Range("F7").Value = 2160
$ws4.Range("F10").Value = 926
$ws4.Range("F12").Value = 998
$ws4.Range("F13").Value = 238
$ws4.Range("F14").Value = 458
$ws4.Range("F19").Value = 7721
$ws4.Range("F20").Value = 337
$ws4.Range("F21").Value = 2454
$ws4.Range("F24").Value = 225
$ws4.Range("F27").Value = 530
$ws4.Range("F30").Value = 977
$ws4.Range("F32").Value = 1649
$ws4.Range("F34").Value = 1165
$ws4.Range("F38").Value = 164
$ws4.Range("F41").Value = 162
$ws4.Range("F42").Value = 333
$ws4.Range("F47").Value = 213
